$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value
$updates = @(
    @{Cell="D2"; Value='29.300.79'},
    @{Cell="E2"; Value='  +0.36%  '},
    @{Cell="D3"; Value='1.874.94'},
    @{Cell="E3"; Value='  +0.64%  '},
    @{Cell="D4"; Value='0.9996'},
    @{Cell="E4"; Value='  -0.24%  '},
    @{Cell="D5"; Value='0.7117'},
    @{Cell="E5"; Value='  +0.25%  '},
    @{Cell="D6"; Value='242.57'},
    @{Cell="E6"; Value='  +0.85%  '},
    @{Cell="E7"; Value='  -0.17%  '},
    @{Cell="E8"; Value='  +1.23%  '},
    @{Cell="D9"; Value='0.07748'},
    @{Cell="E9"; Value='  +0.80%  '},
    @{Cell="D10"; Value='25.06'},
    @{Cell="E10"; Value='  +0.56%  '},
    @{Cell="D11"; Value='0.08459'},
    @{Cell="E11"; Value='  +2.43%  '},
    @{Cell="D12"; Value='1.879.86'},
    @{Cell="E12"; Value='  +1.34%  '},
    @{Cell="E13"; Value='  -0.09%  '},
    @{Cell="D14"; Value='0.7115'},
    @{Cell="E14"; Value='  -0.64%  '},
    @{Cell="D15"; Value='91.29'},
    @{Cell="E15"; Value='  +1.21%  '},
    @{Cell="D16"; Value='29.299.88'},
    @{Cell="E16"; Value='  +0.37%  '},
    @{Cell="D17"; Value='0.000008312'},
    @{Cell="E17"; Value='  +6.70%  '},
    @{Cell="D18"; Value='5.991'},
    @{Cell="E18"; Value='  +2.41%  '},
    @{Cell="D19"; Value='242.52'},
    @{Cell="E19"; Value='  -0.27%  '},
    @{Cell="D20"; Value='13.22'},
    @{Cell="E20"; Value='  +0.65%  '},
    @{Cell="D21"; Value='2.122.02'},
    @{Cell="E21"; Value='  +0.62%  '},
    @{Cell="D22"; Value='0.9994'},
    @{Cell="E22"; Value='  -0.16%  '},
    @{Cell="E23"; Value='  -1.68%  '},
    @{Cell="D24"; Value='0.9999'},
    @{Cell="E24"; Value='  -0.27%  '},
    @{Cell="D25"; Value='0.1610'},
    @{Cell="E25"; Value='  +2.30%  '},
    @{Cell="D26"; Value='162.96'},
    @{Cell="E26"; Value='  +0.33%  '},
    @{Cell="D27"; Value='9.019'},
    @{Cell="E27"; Value='  +1.38%  '},
    @{Cell="E28"; Value='  +1.72%  '},
    @{Cell="D29"; Value='1.515'},
    @{Cell="E29"; Value='  +1.42%  '},
    @{Cell="E30"; Value='  +1.40%  '},
    @{Cell="D31"; Value='4.320'},
    @{Cell="E31"; Value='  +5.58%  '},
    @{Cell="D32"; Value='1.258'},
    @{Cell="E32"; Value='  -4.99%  '},
    @{Cell="D33"; Value='0.05258'},
    @{Cell="E33"; Value='  +1.45%  '},
    @{Cell="E34"; Value='  +1.13%  '},
    @{Cell="D35"; Value='1.174'},
    @{Cell="E35"; Value='  -0.03%  '},
    @{Cell="D36"; Value='0.7446'},
    @{Cell="E36"; Value='  +2.40%  '},
    @{Cell="D37"; Value='2.682'},
    @{Cell="E37"; Value='  -0.11%  '},
    @{Cell="D38"; Value='0.01861'},
    @{Cell="E38"; Value='  +0.83%  '},
    @{Cell="D39"; Value='2.715'},
    @{Cell="E39"; Value='  +1.00%  '},
    @{Cell="D40"; Value='1.169.81'},
    @{Cell="E40"; Value='  +2.49%  '},
    @{Cell="D41"; Value='6.366'},
    @{Cell="E41"; Value='  +4.65%  '},
    @{Cell="D42"; Value='73.03'},
    @{Cell="E42"; Value='  +1.25%  '},
    @{Cell="D43"; Value='0.8866'},
    @{Cell="E43"; Value='  -1.31%  '},
    @{Cell="D44"; Value='106.51'},
    @{Cell="E44"; Value='  +4.99%  '},
    @{Cell="D45"; Value='0.9994'},
    @{Cell="E45"; Value='  -0.17%  '},
    @{Cell="D46"; Value='2.018.93'},
    @{Cell="E46"; Value='  +0.59%  '},
    @{Cell="D47"; Value='1.813'},
    @{Cell="E47"; Value='  +2.77%  '},
    @{Cell="D48"; Value='0.5200'},
    @{Cell="E48"; Value='  -1.30%  '},
    @{Cell="E49"; Value='  +0.40%  '},
    @{Cell="D50"; Value='9.379'},
    @{Cell="E50"; Value='  +1.00%  '},
    @{Cell="D51"; Value='0.4303'},
    @{Cell="E51"; Value='  +1.38%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"   # force text so numeric-looking strings are not coerced to numbers
    $rng.Value = $u.Value
}
